$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 82.98768099999999
$ws.Range("H2").Value = 248.963043
$ws.Range("I2").Value = 0.4489504115427952
$ws.Range("J2").Value = 0.4489504115427952
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.01650666666666667
$ws.Range("N2").Value = 0.04952
$ws.Range("O2").Value = 0.795859985214233
$ws.Range("P2").Value = 0.795859985214233
$ws.Range("Q2").Value = 1.369849987706667
$ws.Range("R2").Value = 12.32864988936
$ws.Range("S2").Value = 0.3573016678923728
$ws.Range("T2").Value = 0.3573016678923728

# Row 3
$ws.Range("G3").Value = 82.98768099999999
$ws.Range("H3").Value = 248.963043
$ws.Range("I3").Value = 0.4489504115427952
$ws.Range("J3").Value = 0.4489504115427952
$ws.Range("O3").Value = 0.2041400147857671
$ws.Range("P3").Value = 0.2041400147857671
$ws.Range("Q3").Value = 0.351369841354
$ws.Range("R3").Value = 3.162328572186
$ws.Range("S3").Value = 0.09164874365042244
$ws.Range("T3").Value = 0.09164874365042244

# Row 4
$ws.Range("G4").Value = 63.14058933333333
$ws.Range("I4").Value = 0.3415807409566563
$ws.Range("J4").Value = 0.3415807409566563
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.01650666666666667
$ws.Range("N4").Value = 0.04952
$ws.Range("O4").Value = 0.795859985214233
$ws.Range("P4").Value = 0.795859985214233
$ws.Range("Q4").Value = 1.042240661262222
$ws.Range("R4").Value = 9.38016595136
$ws.Range("S4").Value = 0.2718504434472312
$ws.Range("T4").Value = 0.2718504434472312

# Row 5
$ws.Range("G5").Value = 63.14058933333333
$ws.Range("I5").Value = 0.3415807409566563
$ws.Range("J5").Value = 0.3415807409566563
$ws.Range("O5").Value = 0.2041400147857671
$ws.Range("P5").Value = 0.2041400147857671
$ws.Range("Q5").Value = 0.2673372552373333
$ws.Range("S5").Value = 0.0697302975094251
$ws.Range("T5").Value = 0.0697302975094251

# Row 6
$ws.Range("I6").Value = 0.2094688475005485
$ws.Range("J6").Value = 0.2094688475005485
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.01650666666666667
$ws.Range("N6").Value = 0.04952
$ws.Range("O6").Value = 0.795859985214233
$ws.Range("P6").Value = 0.795859985214233
$ws.Range("Q6").Value = 0.6391371759466667
$ws.Range("R6").Value = 5.752234583520001
$ws.Range("S6").Value = 0.1667078738746289
$ws.Range("T6").Value = 0.1667078738746289

# Row 7
$ws.Range("I7").Value = 0.2094688475005485
$ws.Range("J7").Value = 0.2094688475005485
$ws.Range("O7").Value = 0.2041400147857671
$ws.Range("P7").Value = 0.2041400147857671
$ws.Range("S7").Value = 0.04276097362591956
$ws.Range("T7").Value = 0.04276097362591956
